$d = $word.ActiveDocument

$pairs = @(
    ,@("36+10=", "20+12=")
    ,@("66-24=", "57-34=")
    ,@("47-4=", "72+9=")
    ,@("57+14=", "89-22=")
    ,@("85-70=", "50-8=")
    ,@("62-9=", "18+43=")
    ,@("69-9=", "17+50=")
    ,@("12+17=", "55-20=")
    ,@("94-33=", "26+44=")
    ,@("40+6=", "30+41=")
    ,@("5+1=", "14+36=")
    ,@("16-1=", "77-6=")
    ,@("29-0=", "38+8=")
    ,@("80-40=", "48-2=")
    ,@("31+64=", "83-61=")
    ,@("63-55=", "20-9=")
    ,@("62+4=", "9+69=")
    ,@("49+26=", "4+37=")
    ,@("58+19=", "10+51=")
    ,@("45-27=", "62+35=")
    ,@("80-27=", "77-66=")
    ,@("12+42=", "34-30=")
    ,@("36-17=", "34-3=")
    ,@("52+47=", "6+28=")
    ,@("70-30=", "90+1=")
    ,@("71-32=", "72-57=")
    ,@("18-11=", "9+81=")
    ,@("68-65=", "86-66=")
    ,@("43+22=", "85-67=")
    ,@("13+49=", "45+26=")
    ,@("61-49=", "67-56=")
    ,@("57+2=", "46-38=")
    ,@("69+1=", "52+46=")
    ,@("74+24=", "32+7=")
    ,@("11-9=", "92-60=")
    ,@("71-17=", "10+57=")
    ,@("38+27=", "83-32=")
    ,@("73-31=", "68-52=")
    ,@("70+2=", "48+32=")
    ,@("84-71=", "65-41=")
    ,@("39-39=", "86+10=")
    ,@("26+38=", "98-69=")
    ,@("56-54=", "32+32=")
    ,@("9+22=", "33+42=")
    ,@("17+55=", "32-22=")
    ,@("17-9=", "22+19=")
    ,@("49-23=", "80-77=")
    ,@("8+79=", "83-39=")
    ,@("82-8=", "2+85=")
    ,@("44+1=", "46+38=")
    ,@("13-10=", "27-5=")
    ,@("59+13=", "64-51=")
    ,@("23+47=", "55+24=")
    ,@("92-37=", "67-2=")
    ,@("39-16=", "93-76=")
    ,@("73-61=", "54-46=")
    ,@("58+4=", "64-48=")
    ,@("86-77=", "0+44=")
    ,@("89-6=", "16+23=")
    ,@("63-40=", "70-17=")
    ,@("34-11=", "54+1=")
    ,@("42+6=", "72-8=")
    ,@("87-17=", "3+7=")
    ,@("66-12=", "78-37=")
    ,@("47-46=", "14+63=")
    ,@("68+11=", "26+68=")
    ,@("40-18=", "4+14=")
    ,@("94+5=", "47+12=")
    ,@("79-44=", "40+47=")
    ,@("43+48=", "3+40=")
    ,@("58+15=", "44+28=")
    ,@("25+26=", "98-54=")
    ,@("80-54=", "87-2=")
    ,@("22+33=", "13+45=")
    ,@("52-20=", "30+56=")
    ,@("86-51=", "59-16=")
    ,@("50-20=", "5+32=")
    ,@("11+28=", "94-23=")
    ,@("86-42=", "73-18=")
    ,@("51-34=", "52+20=")
    ,@("24+18=", "85-32=")
    ,@("28+51=", "14+64=")
    ,@("21+12=", "21-9=")
    ,@("6+45=", "1+47=")
    ,@("54-23=", "33+38=")
    ,@("32+28=", "84-37=")
    ,@("18-8=", "74-21=")
    ,@("43-19=", "94-13=")
    ,@("88-43=", "35+18=")
    ,@("3+9=", "6+27=")
    ,@("69-48=", "33+55=")
    ,@("79+4=", "10+73=")
    ,@("84-15=", "37-22=")
    ,@("11+30=", "13+6=")
    ,@("88+8=", "60+23=")
    ,@("30+63=", "30+1=")
    ,@("95-35=", "0+96=")
    ,@("18+9=", "51+42=")
    ,@("10+68=", "68+31=")
    ,@("45+24=", "97-30=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
